# Updated symbol list on Fri Jan 27 22:10:43 UTC 2023 with GitHub Actions
# Refresh Price / Volume(1h) / Hora columns for the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number plus the new Price (D), Volume(1h) (E) and Hora (G) text.
# Missing keys mean that column is unchanged for that row.
$updates = @(
  @{ Row=2; D="307.32"; E="0.62%"; G="22" },
  @{ Row=3; D="36.47"; E="1.58%"; G="22" },
  @{ Row=4; D="5.067"; E="1.54%"; G="22" },
  @{ Row=5; D="0.08095"; E="-0.08%"; G="22" },
  @{ Row=6; D="2.011"; E="5.32%"; G="22" },
  @{ Row=7; D="7.856"; E="-0.29%"; G="22" },
  @{ Row=8; D="0.9282"; E="-0.19%"; G="22" },
  @{ Row=9; D="0.1468"; E="13.86%"; G="22" },
  @{ Row=10; D="0.1945"; E="2.29%"; G="22" },
  @{ Row=11; D="0.09145"; E="-0.68%"; G="22" },
  @{ Row=12; D="0.03519"; E="0.09%"; G="22" },
  @{ Row=13; D="0.09888"; E="-0.23%"; G="22" },
  @{ Row=14; D="0.001407"; E="-1.66%"; G="22" },
  @{ Row=15; D="0.006388"; E="-3.45%"; G="22" },
  @{ Row=16; D="3.840"; E="6.53%"; G="22" },
  @{ Row=17; D="4.173"; E="0.68%"; G="22" },
  @{ Row=18; E="7.48%"; G="22" },
  @{ Row=19; D="0.3451"; E="-0.08%"; G="22" },
  @{ Row=20; D="0.1327"; E="2.38%"; G="22" },
  @{ Row=21; D="4.834"; E="-7.97%"; G="22" },
  @{ Row=22; D="0.2348"; E="-7.22%"; G="22" },
  @{ Row=23; D="0.04399"; E="-0.24%"; G="22" },
  @{ Row=24; E="0.06%"; G="22" },
  @{ Row=25; D="0.004181"; E="-11.23%"; G="22" },
  @{ Row=26; G="22" },
  @{ Row=27; D="0.0001304"; E="0.49%"; G="22" },
  @{ Row=28; G="22" },
  @{ Row=29; G="22" },
  @{ Row=30; G="22" },
  @{ Row=31; G="22" },
  @{ Row=32; G="22" },
  @{ Row=33; G="22" },
  @{ Row=34; G="22" },
  @{ Row=35; G="22" },
  @{ Row=36; G="22" },
  @{ Row=37; G="22" },
  @{ Row=38; G="22" },
  @{ Row=39; D="0.02046"; E="4.84%"; G="22" },
  @{ Row=40; D="0.05131"; E="-1.60%"; G="22" },
  @{ Row=41; D="0.007479"; E="-1.06%"; G="22" },
  @{ Row=42; D="0.01009"; E="-0.96%"; G="22" },
  @{ Row=43; D="0.1369"; E="-0.10%"; G="22" },
  @{ Row=44; D="0.002127"; E="1.46%"; G="22" },
  @{ Row=45; D="0.009899"; E="-6.97%"; G="22" },
  @{ Row=46; D="0.00006316"; E="-0.72%"; G="22" },
  @{ Row=47; D="0.00000000752"; E="0.26%"; G="22" },
  @{ Row=48; E="-0.16%"; G="22" },
  @{ Row=49; D="0.001605"; E="-3.29%"; G="22" },
  @{ Row=50; D="0.00002107"; E="0.26%"; G="22" },
  @{ Row=51; D="0.0002007"; E="0.26%"; G="22" }
)

foreach ($item in $updates) {
  $r = $item.Row
  foreach ($col in @("D", "E", "G")) {
    if ($item.ContainsKey($col)) {
      $cell = $ws.Range("$col$r")
      # Force text storage ("@") so numeric-looking strings (e.g. "307.32", "22")
      # and percentages (e.g. "0.62%") are kept as literal text, matching the
      # original inline-string cell type instead of being parsed into numbers.
      $cell.NumberFormat = "@"
      $cell.Value = $item[$col]
    }
  }
}